# Edit: Added Negative Login Test Cases + Added XFail pytest marks
#
# Summary of changes (per commit diff):
# 1. "Login - Negative" sheet: remove the brute-force (old TC-LOGIN-NEG-001) and
#    invalid-username-format (old TC-LOGIN-NEG-003) test rows entirely.
# 2. Renumber the remaining Test Case IDs on that sheet (NEG-002->001, NEG-004->002,
#    NEG-005->003, NEG-006->004, NEG-007->005; NEG-008 stays NEG-008).
# 3. Extend "Login - Negative" sheet/table with the Timestamp / Duration (s) /
#    Actual Result columns (mirroring the "Login - Positive" sheet) and resize table.
# 4. Make "Login - Negative" the active/selected sheet & tab; update each sheet's
#    remembered selection cell.

$wb = $excel.ActiveWorkbook
$wsPos = $wb.Worksheets.Item("Login - Positive")
$wsNeg = $wb.Worksheets.Item("Login - Negative")

# --- 1 & 2: delete obsolete rows, renumber remaining Test Case IDs -----------------
# Current (before) layout of "Login - Negative":
#   row2 TC-LOGIN-NEG-001 Try multiple invalid logins (brute-force)        -> delete
#   row3 TC-LOGIN-NEG-002 Enter mixed-case username ...                    -> becomes NEG-001
#   row4 TC-LOGIN-NEG-003 Enter invalid username format ...                -> delete
#   row5 TC-LOGIN-NEG-004 Enter unregistered username ...                  -> becomes NEG-002
#   row6 TC-LOGIN-NEG-005 Enter valid username with incorrect password...  -> becomes NEG-003
#   row7 TC-LOGIN-NEG-006 Leave username field empty ...                   -> becomes NEG-004
#   row8 TC-LOGIN-NEG-007 Enter valid username, leave password field empty -> becomes NEG-005
#   row9 TC-LOGIN-NEG-008 Leave both fields empty ...                      -> stays NEG-008

$wsNeg.Rows.Item(4).Delete()
$wsNeg.Rows.Item(2).Delete()

$wsNeg.Range("A2").Value = "TC-LOGIN-NEG-001"
$wsNeg.Range("A3").Value = "TC-LOGIN-NEG-002"
$wsNeg.Range("A4").Value = "TC-LOGIN-NEG-003"
$wsNeg.Range("A5").Value = "TC-LOGIN-NEG-004"
$wsNeg.Range("A6").Value = "TC-LOGIN-NEG-005"
$wsNeg.Range("A7").Value = "TC-LOGIN-NEG-008"

# --- 3: add Timestamp / Duration (s) / Actual Result columns -----------------------
$loNeg = $wsNeg.ListObjects.Item(1)
$loNeg.Resize($wsNeg.Range("A1:I7"))
$wsNeg.Range("G1").Value = "Timestamp"
$wsNeg.Range("H1").Value = "Duration (s)"
$wsNeg.Range("I1").Value = "Actual Result"

# --- 4: active sheet / selections ---------------------------------------------------
$wsNeg.Activate()
$wsNeg.Range("B12").Select()
$wsPos.Range("D12").Select()
$wsNeg.Activate()
